$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-05 Tuesday" "2023-12-06 Wednesday"

Replace-Text "17×51=" "50×40="
Replace-Text "74×49=" "43×19="
Replace-Text "68×25=" "73×15="
Replace-Text "55×45=" "96×93="
Replace-Text "63×56=" "91×90="
Replace-Text "31×56=" "98×67="
Replace-Text "13×48=" "36×36="
Replace-Text "88×79=" "87×48="
Replace-Text "98×17=" "60×79="
Replace-Text "55×40=" "16×61="
Replace-Text "47×83=" "68×13="
Replace-Text "75×56=" "29×69="
Replace-Text "92×60=" "50×94="
Replace-Text "44×45=" "79×77="
Replace-Text "97×25=" "59×75="
Replace-Text "68×51=" "67×51="
Replace-Text "45×49=" "94×91="
Replace-Text "37×98=" "93×15="
Replace-Text "44×79=" "70×79="
Replace-Text "89×39=" "93×64="
Replace-Text "70×73=" "30×27="
Replace-Text "75×24=" "69×66="
Replace-Text "39×44=" "43×71="
Replace-Text "21×12=" "65×18="
Replace-Text "99×76=" "64×75="
